# Update rules in DiscountRules.xlsx
# - Row 24: set A24/B24/C24 to "Test"
# - Row 25: set B25/C25 to "Test" (A25 already "Test")
# - Remove row 26 entirely (shifts dimension from A1:C26 to A1:C25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "Test"
$ws.Range("B24").Value = "Test"
$ws.Range("C24").Value = "Test"

$ws.Range("B25").Value = "Test"
$ws.Range("C25").Value = "Test"

$ws.Rows("26").Delete()
